# Update "想去人数" (F) and "最低票价" (G) figures on the
# "展览" and "全部类型" sheets to match the newly scraped data.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F2").Value = 1083
    $ws.Range("G2").Value = 65

    $ws.Range("F7").Value = 2422

    $ws.Range("F11").Value = 1214

    $ws.Range("F15").Value = 1084

    $ws.Range("F16").Value = 298

    $ws.Range("F18").Value = 20

    $ws.Range("F23").Value = 114

    $ws.Range("F24").Value = 10

    $ws.Range("F25").Value = 243
}
